$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Widen column J (column index 10) from stored width 9 to 11.
# The COM ColumnWidth property is offset by ~0.83 from the raw
# OOXML "width" attribute for this workbook's default font, so we
# subtract that offset to land on the target stored width of 11.
$ws1.Columns.Item(10).ColumnWidth = 10.17

# Row 4 (FARIAS CAICEDO GABRIELA PATRICIA)
$ws1.Range("D4").Value = 434.83

# Row 7 (MOROCHO PLAZA SHIRLEY AURELIA)
$ws1.Range("J7").Value = 36.74
$ws1.Range("N7").Value = 100.71
$ws1.Range("Q7").Value = 21.58

# Row 8 ("x de 6" completion counters)
$ws1.Range("D8").Value = "2 de 6"
$ws1.Range("J8").Value = "1 de 6"
$ws1.Range("N8").Value = "1 de 6"
$ws1.Range("Q8").Value = "1 de 6"

# ---------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 4 (FARIAS CAICEDO GABRIELA PATRICIA)
$ws2.Range("F4").Value = 434.83

# Row 7 (MOROCHO PLAZA SHIRLEY AURELIA)
$ws2.Range("F7").Value = 159.03

# Row 8 (total)
$ws2.Range("F8").Value = 1066.43

# ---------------------------------------------------------------
# Sheet 3: "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS group)
$ws3.Range("D2").Value = 5043.87
$ws3.Range("E2").Value = -5043.87

# Row 4 (TOTAL)
$ws3.Range("D4").Value = 5638.65
$ws3.Range("E4").Value = 11861.35
$ws3.Range("F4").Value = 0.3222085714285714
